$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.965.59"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "3.103.14"
$ws.Range("E3").Value = "  +5.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'580.05"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "'173.98"
$ws.Range("E6").Value = "  +8.22%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.097.05"
$ws.Range("E8").Value = "  +5.26%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "'6.51"
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("E11").Value = "  +4.44%  "
$ws.Range("D12").Value = "'0.482"
$ws.Range("E12").Value = "  +5.35%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").Value = "'37.26"
$ws.Range("E14").Value = "  +8.02%  "
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "3.614.94"
$ws.Range("E16").Value = "  +6.05%  "
$ws.Range("D17").Value = "66.968.54"
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").Value = "'7.20"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").Value = "3.103.11"
$ws.Range("E19").Value = "  +5.36%  "
$ws.Range("D20").Value = "'16.19"
$ws.Range("E20").Value = "  +3.22%  "
$ws.Range("D21").Value = "'481.34"
$ws.Range("E21").Value = "  +8.22%  "
$ws.Range("D22").Value = "'0.716"
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("D24").Value = "'84.17"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +5.38%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'13.04"
$ws.Range("E26").Value = "  +7.22%  "
$ws.Range("D27").Value = "'10.01"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D29").Value = "'8.02"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").Value = "'2.39"
$ws.Range("E30").Value = "  -3.02%  "
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").Value = "'28.77"
$ws.Range("E33").Value = "  +5.88%  "
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").Value = "'5.90"
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("D38").Value = "'47.67"
$ws.Range("E38").Value = "  +5.61%  "
$ws.Range("D39").Value = "'2.13"
$ws.Range("E39").Value = "  +7.10%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.317"
$ws.Range("E40").Value = "  +5.61%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'50.17"
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").Value = "'8.68"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "'0.0361"
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("D46").Value = "2.821.29"
$ws.Range("E46").Value = "  +5.41%  "
$ws.Range("D47").Value = "'380.90"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'134.96"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "'24.90"
$ws.Range("E50").Value = "  +5.48%  "
$ws.Range("E51").Value = "  +2.46%  "
